$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the header row ("fase_dia") - row 1. Remaining rows shift up.
$ws.Rows.Item(1).Delete()

# Insert a new row before the current row 3 ("Plena Noite") and set it to "n/a".
$ws.Rows.Item(3).Insert()
$ws.Range("A3").Value = "n/a"

# Update the selection to the full column A (mirrors the post-edit selection in the file).
$ws.Range("A:A").Select()
